$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 82 (existing rows 82-139 shift down to 83-140)
$ws.Rows.Item(82).Insert()

# Populate the new row 82 with the data for the added price record
$ws.Cells.Item(82, 1).Value = 5
$ws.Cells.Item(82, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(82, 3).Value = "Maule"
$ws.Cells.Item(82, 4).Value = 44438
$ws.Cells.Item(82, 5).Value = 7
$ws.Cells.Item(82, 6).Value = "Fruta"
$ws.Cells.Item(82, 7).Value = 100108
$ws.Cells.Item(82, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(82, 9).Value = 100108005
$ws.Cells.Item(82, 10).Value = "Piña"
$ws.Cells.Item(82, 11).Value = "Caramelo"
$ws.Cells.Item(82, 12).Value = "Tercera"
$ws.Cells.Item(82, 13).Value = 200
$ws.Cells.Item(82, 14).Value = 20000
$ws.Cells.Item(82, 15).Value = 20000
$ws.Cells.Item(82, 16).Value = 20000
$ws.Cells.Item(82, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(82, 18).Value = "Ecuador"
$ws.Cells.Item(82, 19).Value = 1250
$ws.Cells.Item(82, 20).Value = 16
